$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the departure/return timestamps in H6/I6 and H8/I8 -
# shift the date portion back one day (Nov 2 -> Nov 1), keep time-of-day.
$ws.Range("H6").Value = 44866.333333333336
$ws.Range("I6").Value = 44866.708333333336
$ws.Range("H8").Value = 44866.333333333336
$ws.Range("I8").Value = 44866.708333333336

# Update the active selection shown when the file is opened.
$ws.Range("G32").Select()
